# Weekly "Hortaliza, Vega Monumental Concepción - Acelga" refresh.
# A new weekly reading (rows 129-130, date 2021-12-23) is inserted at the top of
# the "Acelga" block, which pushes the dates for every later fortnightly pair down
# by one slot; two more rows (193-194) are appended at the bottom to restore the
# two entries that fall off the end of the shifted sequence. A couple of older rows
# also get small data corrections (volumes, region labels, unit label).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new top-of-block reading: 2021-12-23 (row 129/130), unit label simplified ---
$ws.Range("D129").Value = 44553
$ws.Range("N129").Value = "`$/atado"
$ws.Range("D130").Value = 44553
$ws.Range("N130").Value = "`$/atado"

# --- volumes corrected on the next two fortnightly pairs ---
$ws.Range("D131").Value = 44161
$ws.Range("J131").Value = 200
$ws.Range("D132").Value = 44161
$ws.Range("J132").Value = 100
$ws.Range("D133").Value = 44285
$ws.Range("J133").Value = 100
$ws.Range("D134").Value = 44285
$ws.Range("J134").Value = 50

# --- remaining pairs: date shifts down the sequence by one slot ---
$ws.Range("D135").Value = 44434
$ws.Range("D136").Value = 44434
$ws.Range("D137").Value = 44467
$ws.Range("D138").Value = 44467
$ws.Range("D139").Value = 44231
$ws.Range("D140").Value = 44231
$ws.Range("D141").Value = 44490
$ws.Range("D142").Value = 44490
$ws.Range("D143").Value = 44259
$ws.Range("D144").Value = 44259
$ws.Range("D145").Value = 44341
$ws.Range("D146").Value = 44341
$ws.Range("D147").Value = 44208
$ws.Range("D148").Value = 44208
$ws.Range("D149").Value = 44264
$ws.Range("D150").Value = 44264
$ws.Range("D151").Value = 44322
$ws.Range("D152").Value = 44322
$ws.Range("D153").Value = 44391
$ws.Range("D154").Value = 44391
$ws.Range("D155").Value = 44551
$ws.Range("D156").Value = 44551
$ws.Range("D157").Value = 44386
$ws.Range("D158").Value = 44386
$ws.Range("D159").Value = 44420
$ws.Range("D160").Value = 44420
$ws.Range("D161").Value = 44278
$ws.Range("D162").Value = 44278
$ws.Range("D163").Value = 44308
$ws.Range("D164").Value = 44308
$ws.Range("D165").Value = 44187
$ws.Range("D166").Value = 44187
$ws.Range("D167").Value = 44474
$ws.Range("D168").Value = 44474
$ws.Range("D169").Value = 44350
$ws.Range("D170").Value = 44350
$ws.Range("D171").Value = 44405
$ws.Range("D172").Value = 44405
$ws.Range("D173").Value = 44257
$ws.Range("D174").Value = 44257
$ws.Range("D175").Value = 44196
$ws.Range("D176").Value = 44196
$ws.Range("D177").Value = 44224
$ws.Range("D178").Value = 44224

# --- these two pairs also swap their "Origen" region label ---
$ws.Range("D179").Value = 44398
$ws.Range("O179").Value = "Región de Ñuble"
$ws.Range("D180").Value = 44398
$ws.Range("O180").Value = "Región de Ñuble"
$ws.Range("D181").Value = 44239
$ws.Range("O181").Value = "Región Metropolitana"
$ws.Range("D182").Value = 44239
$ws.Range("O182").Value = "Región Metropolitana"

# --- remaining pairs continue the date shift ---
$ws.Range("D183").Value = 44344
$ws.Range("D184").Value = 44344
$ws.Range("D185").Value = 44371
$ws.Range("D186").Value = 44371
$ws.Range("D187").Value = 44365
$ws.Range("D188").Value = 44365
$ws.Range("D189").Value = 44194
$ws.Range("D190").Value = 44194

# --- row 191 shifts date and is corrected to the "Primera" volumes/prices ---
$ws.Range("D191").Value = 44313
$ws.Range("J191").Value = 200
$ws.Range("L191").Value = 700
$ws.Range("M191").Value = 650
$ws.Range("P191").Value = 650

# --- row 192 becomes the "Segunda" reading that now shares row 191's new date ---
$ws.Range("D192").Value = 44313
$ws.Range("I192").Value = "Segunda"
$ws.Range("J192").Value = 100
$ws.Range("L192").Value = 500
$ws.Range("M192").Value = 500
$ws.Range("P192").Value = 500

# --- two brand-new rows appended: the readings bumped off the end of the block ---
$ws.Range("A193").Value = 11
$ws.Range("B193").Value = "Vega Monumental Concepción"
$ws.Range("C193").Value = "Bíobío"
$ws.Range("D193").Value = 44518
$ws.Range("D193").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E193").Value = 8
$ws.Range("F193").Value = 100112009
$ws.Range("G193").Value = "Acelga"
$ws.Range("H193").Value = "Sin especificar"
$ws.Range("I193").Value = "Primera"
$ws.Range("J193").Value = 450
$ws.Range("K193").Value = 600
$ws.Range("L193").Value = 650
$ws.Range("M193").Value = 628
$ws.Range("N193").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O193").Value = "Región de Ñuble"
$ws.Range("P193").Value = 628
$ws.Range("Q193").Value = 1
$ws.Range("R193").Value = "Hortaliza"

$ws.Range("A194").Value = 11
$ws.Range("B194").Value = "Vega Monumental Concepción"
$ws.Range("C194").Value = "Bíobío"
$ws.Range("D194").Value = 44540
$ws.Range("D194").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E194").Value = 8
$ws.Range("F194").Value = 100112009
$ws.Range("G194").Value = "Acelga"
$ws.Range("H194").Value = "Sin especificar"
$ws.Range("I194").Value = "Primera"
$ws.Range("J194").Value = 450
$ws.Range("K194").Value = 500
$ws.Range("L194").Value = 550
$ws.Range("M194").Value = 522
$ws.Range("N194").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O194").Value = "Región de Ñuble"
$ws.Range("P194").Value = 522
$ws.Range("Q194").Value = 1
$ws.Range("R194").Value = "Hortaliza"
